# Generate Report for Handoff
# Adds two new tracked files (a .png pair and a .md file) to the
# localization-status workbook: updates the existing "1876bef4...png"
# (formerly "a31496fd...md") row with its refreshed handoff info, and
# appends two new rows (one per new source file) across all three sheets:
#   Overview (sheet 1), zh-cn (sheet 2), de-de (sheet 3)

function Set-LinkCell {
    param(
        $ws,
        [string]$cellRef,
        [string]$displayText,
        [string]$targetUrl
    )
    $rng = $ws.Range($cellRef)
    # Drop any pre-existing hyperlink bound to this cell first so we don't
    # leave a stale/duplicate entry behind, then re-add cleanly.
    $rng.Hyperlinks.Delete()
    $rng.Value = $displayText
    $ws.Hyperlinks.Add($rng, $targetUrl, "", "", $displayText) | Out-Null
}

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

# ---- source / artifact identifiers ------------------------------------
$png1        = "1876bef4-39ab-449e-bf32-777224809fd8.png"
$png2        = "d68ea10c-70b1-4507-b40d-278a9fdfb4f7.png"
$mdFile      = "f7704d6c-4468-42c1-a02c-d4a42962623a.md"

$png1Target  = "c809d7b7da740dfacc6f4c7016a657580d75b655.png"
$png2Target  = "d960e153e778bf8393efb403f328f7d27cd7a9b1.png"
$mdZhTarget  = "f7704d6c-4468-42c1-a02c-d4a42962623a.7b18c1365f523eb400628e33aef23dfe4b0ffe4a.zh-cn.xlf"
$mdDeTarget  = "f7704d6c-4468-42c1-a02c-d4a42962623a.7b18c1365f523eb400628e33aef23dfe4b0ffe4a.de-de.xlf"

$dependencyFrom = "e2e\f7704d6c-4468-42c1-a02c-d4a42962623a.md"

$readyStatus = "Ready for handoff"

$overviewDate = "2016-49-12 16:49:16"
$zhDate       = "2016-03-12 16:49:13"
$deDate       = "2016-03-12 16:49:16"
$epoch        = "0001-01-01 00:00:00"

# github source-file blob urls (column A / B)
$srcUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/c896770c956a536cd9c39e8f254743774594b8a4/e2e/$png1"
$srcUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/c896770c956a536cd9c39e8f254743774594b8a4/e2e/$png2"
$srcUrl3 = "https://github.com/OpenLocalizationTest/oltest/blob/c896770c956a536cd9c39e8f254743774594b8a4/e2e/$mdFile"

# olhandoff target-file blob urls (column D)
$tgtUrl1zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb76e0c7ab1e1c6ac017312c5346fe0a2c3af2b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1Target"
$tgtUrl2zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb76e0c7ab1e1c6ac017312c5346fe0a2c3af2b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2Target"
$tgtUrl3zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb76e0c7ab1e1c6ac017312c5346fe0a2c3af2b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$mdZhTarget"

$tgtUrl1de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a58e804df85be97d4d1fea26d66d8de875772a4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1Target"
$tgtUrl2de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a58e804df85be97d4d1fea26d66d8de875772a4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2Target"
$tgtUrl3de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a58e804df85be97d4d1fea26d66d8de875772a4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$mdDeTarget"

# =========================================================================
# Overview sheet: columns A=File Name B=zh-cn C=de-de D=Latest Handoff Date
# =========================================================================

# existing row 2 -> now describes the first .png file
Set-LinkCell $overview "A2" $png1 $srcUrl1
$overview.Range("B2").Value = $readyStatus
$overview.Range("C2").Value = $readyStatus
$overview.Range("D2").Value = $overviewDate

# new row 3 -> second .png file
Set-LinkCell $overview "A3" $png2 $srcUrl2
$overview.Range("B3").Value = $readyStatus
$overview.Range("C3").Value = $readyStatus
$overview.Range("D3").Value = $overviewDate

# new row 4 -> the .md file
Set-LinkCell $overview "A4" $mdFile $srcUrl3
$overview.Range("B4").Value = $readyStatus
$overview.Range("C4").Value = $readyStatus
$overview.Range("D4").Value = $overviewDate

# =========================================================================
# zh-cn / de-de sheets: columns
# A=Source File Name B=File Extension C=Status D=Latest Handoff File
# E=Latest Handoff Datetime F=Latest Target File G=Latest Handback File
# H=Latest Handback DateTime I=Handoff Reason J=Dependency From K=Error Detail
# =========================================================================

foreach ($pair in @(
        @{ ws = $zhcn; date = $zhDate; tgt1 = $tgtUrl1zh; tgt2 = $tgtUrl2zh; tgt3 = $tgtUrl3zh; mdTarget = $mdZhTarget },
        @{ ws = $dede; date = $deDate; tgt1 = $tgtUrl1de; tgt2 = $tgtUrl2de; tgt3 = $tgtUrl3de; mdTarget = $mdDeTarget }
    )) {

    $ws = $pair.ws

    # ---- row 2: first .png file (IsDependency) -------------------------
    Set-LinkCell $ws "A2" $png1 $srcUrl1
    Set-LinkCell $ws "B2" ".png" $srcUrl1
    $ws.Range("C2").Value = $readyStatus
    Set-LinkCell $ws "D2" $png1Target $pair.tgt1
    $ws.Range("E2").Value = $pair.date
    $ws.Range("H2").Value = $epoch
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = $dependencyFrom

    # ---- row 3: second .png file (IsDependency) ------------------------
    Set-LinkCell $ws "A3" $png2 $srcUrl2
    Set-LinkCell $ws "B3" ".png" $srcUrl2
    $ws.Range("C3").Value = $readyStatus
    Set-LinkCell $ws "D3" $png2Target $pair.tgt2
    $ws.Range("E3").Value = $pair.date
    $ws.Range("H3").Value = $epoch
    $ws.Range("I3").Value = "IsDependency"
    $ws.Range("J3").Value = $dependencyFrom

    # ---- row 4: the .md file (Include) ---------------------------------
    Set-LinkCell $ws "A4" $mdFile $srcUrl3
    Set-LinkCell $ws "B4" ".md" $srcUrl3
    $ws.Range("C4").Value = $readyStatus
    Set-LinkCell $ws "D4" $pair.mdTarget $pair.tgt3
    $ws.Range("E4").Value = $pair.date
    $ws.Range("H4").Value = $epoch
    $ws.Range("I4").Value = "Include"
}
